$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), reusing the same
# formatting (bold/centered/bordered) as the existing header cells by
# copying an existing header cell's format+value first, then overwriting
# the value/text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I (I0) and J (IF)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 5

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 6

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 5

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 6
